$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New predictions for 2021-01-09, generated after incorporating weather
# data into the model. Two new weekly rows are appended to the table:
#   row 52: week "10 Jan -- 16 Jan 2021" (same week as row 50) with the
#           weather-updated Weekly MAE / Weekly MAPE figures
#   row 53: week "17 Jan -- 23 Jan 2021" (same week as row 51), prediction
#           only, unchanged from before
#
# Column A holds the value "2021-01-09" which must be stored as literal
# text (matching the rest of the column), not auto-converted to a date
# serial number. Pre-formatting the cell as Text ("@") before assigning
# the value keeps it as a string; resetting the style back to Normal
# afterwards avoids leaving a visible Text format applied to the cell.
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "2021-01-09"
$ws.Range("A52").Style = "Normal"

$ws.Range("B52").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("C52").Value = 3333.57
$ws.Range("D52").Value = 2156.57
$ws.Range("E52").Value = 1177.01
$ws.Range("F52").Value = "KNN"
$ws.Range("J52").Value = 966.35
$ws.Range("K52").Value = 31.6

$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "2021-01-09"
$ws.Range("A53").Style = "Normal"

$ws.Range("B53").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D53").Value = 2493.59
$ws.Range("F53").Value = "KNN"
